# ============================================================================
# "Chiffres  COVID-19 Valais.xlsx" data refresh
# Title: "Donnees COVID-19 Valais 20.05.2020" -> "...25.05.2020"
# Historic "Patients COVID-19 hospitalises hors SI" (col G) series revised
# upward (rows 3-76); a handful of "Nb nouvelles admissions" (col D) values
# revised too; six new daily rows (86-90, 21-25 May) appended, and the old
# "today, not final" placeholder row (85) becomes an ordinary completed row
# while row 90 becomes the new placeholder / "today" row.
# ============================================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the report title held in A1 (merged A1:L1) -----------------
$ws.Range("A1").Value = "Données COVID-19 Valais 25.05.2020"

# --- 2. Revise a few "Nb nouvelles admissions a l'hopital" (col D) values -
$ws.Range("D3").Value = 7
$ws.Range("D59").Value = 3
$ws.Range("D64").Value = 2
$ws.Range("D69").Value = 1

# --- 3. Revise the historical "hospitalises hors SI" (col G) series -------
# (col H = G+E recalculates automatically)
$ws.Range("G3").Value = 7
$ws.Range("G4").Value = 7
$ws.Range("G5").Value = 8
$ws.Range("G6").Value = 11
$ws.Range("G7").Value = 13
$ws.Range("G8").Value = 14
$ws.Range("G9").Value = 12
$ws.Range("G10").Value = 13
$ws.Range("G11").Value = 13
$ws.Range("G12").Value = 16
$ws.Range("G13").Value = 17
$ws.Range("G14").Value = 20
$ws.Range("G15").Value = 21
$ws.Range("G16").Value = 28
$ws.Range("G17").Value = 30
$ws.Range("G18").Value = 34
$ws.Range("G19").Value = 36
$ws.Range("G20").Value = 41
$ws.Range("G21").Value = 46
$ws.Range("G22").Value = 58
$ws.Range("G23").Value = 60
$ws.Range("G24").Value = 68
$ws.Range("G25").Value = 76
$ws.Range("G26").Value = 82
$ws.Range("G27").Value = 93
$ws.Range("G28").Value = 106
$ws.Range("G29").Value = 105
$ws.Range("G30").Value = 111
$ws.Range("G31").Value = 119
$ws.Range("G32").Value = 123
$ws.Range("G33").Value = 132
$ws.Range("G34").Value = 130
$ws.Range("G35").Value = 130
$ws.Range("G36").Value = 122
$ws.Range("G37").Value = 120
$ws.Range("G38").Value = 122
$ws.Range("G39").Value = 123
$ws.Range("G40").Value = 122
$ws.Range("G41").Value = 109
$ws.Range("G42").Value = 107
$ws.Range("G43").Value = 99
$ws.Range("G44").Value = 97
$ws.Range("G45").Value = 95
$ws.Range("G46").Value = 96
$ws.Range("G47").Value = 92
$ws.Range("G48").Value = 89
$ws.Range("G49").Value = 84
$ws.Range("G50").Value = 78
$ws.Range("G51").Value = 79
$ws.Range("G52").Value = 73
$ws.Range("G53").Value = 70
$ws.Range("G54").Value = 70
$ws.Range("G55").Value = 68
$ws.Range("G56").Value = 66
$ws.Range("G57").Value = 68
$ws.Range("G58").Value = 64
$ws.Range("G59").Value = 61
$ws.Range("G60").Value = 63
$ws.Range("G61").Value = 64
$ws.Range("G62").Value = 56
$ws.Range("G63").Value = 53
$ws.Range("G64").Value = 55
$ws.Range("G65").Value = 56
$ws.Range("G66").Value = 54
$ws.Range("G67").Value = 54
$ws.Range("G68").Value = 55
$ws.Range("G69").Value = 51
$ws.Range("G70").Value = 49
$ws.Range("G71").Value = 46
$ws.Range("G72").Value = 41
$ws.Range("G73").Value = 39
$ws.Range("G74").Value = 37
$ws.Range("G75").Value = 38
$ws.Range("G76").Value = 38

# --- 4. Turn the former placeholder row 85 into a normal completed row, ---
#        then append new rows 86-90 (new data through 25.05.2020), and    -
#        finally turn row 90 into the new placeholder / "today" row.      -
#        Formats are copied down before values are written so every new   -
#        row keeps the same look (borders/fill/number-format) as row 84.  -

# Capture the placeholder-row look (row 85) before we overwrite it, so it
# can be re-applied to the new placeholder row (90) afterwards.
$ws.Range("A85:L85").Copy()
$ws.Range("A90:L90").PasteSpecial(-4122)

# Re-format rows 85-89 as normal data rows (same look as row 84).
$ws.Range("A84:L84").Copy()
$ws.Range("A85:L89").PasteSpecial(-4122)

# Row 77 (43963)
$ws.Range("A77").Value = 43963
$ws.Range("B77").Formula = "=B76+C77"
$ws.Range("C77").Value = 0
$ws.Range("D77").Value = 1
$ws.Range("E77").Value = 7
$ws.Range("F77").Value = 6
$ws.Range("G77").Value = 38
$ws.Range("H77").Formula = "=G77+E77"
$ws.Range("I77").Formula = "=I76+J77"
$ws.Range("J77").Formula = "=K77+L77"
$ws.Range("K77").NumberFormat = "General"
$ws.Range("K77").Value = 0
$ws.Range("K77").NumberFormat = "@"
$ws.Range("L77").NumberFormat = "General"
$ws.Range("L77").Value = 0
$ws.Range("L77").NumberFormat = "@"

# Row 78 (43964)
$ws.Range("A78").Value = 43964
$ws.Range("B78").Formula = "=B77+C78"
$ws.Range("C78").Value = 1
$ws.Range("D78").Value = 0
$ws.Range("E78").Value = 7
$ws.Range("F78").Value = 5
$ws.Range("G78").Value = 35
$ws.Range("H78").Formula = "=G78+E78"
$ws.Range("I78").Formula = "=I77+J78"
$ws.Range("J78").Formula = "=K78+L78"
$ws.Range("K78").NumberFormat = "General"
$ws.Range("K78").Value = 0
$ws.Range("K78").NumberFormat = "@"
$ws.Range("L78").NumberFormat = "General"
$ws.Range("L78").Value = 0
$ws.Range("L78").NumberFormat = "@"

# Row 79 (43965)
$ws.Range("A79").Value = 43965
$ws.Range("B79").Formula = "=B78+C79"
$ws.Range("C79").Value = 2
$ws.Range("D79").Value = 0
$ws.Range("E79").Value = 7
$ws.Range("F79").Value = 5
$ws.Range("G79").Value = 33
$ws.Range("H79").Formula = "=G79+E79"
$ws.Range("I79").Formula = "=I78+J79"
$ws.Range("J79").Formula = "=K79+L79"
$ws.Range("K79").NumberFormat = "General"
$ws.Range("K79").Value = 0
$ws.Range("K79").NumberFormat = "@"
$ws.Range("L79").NumberFormat = "General"
$ws.Range("L79").Value = 0
$ws.Range("L79").NumberFormat = "@"

# Row 80 (43966)
$ws.Range("A80").Value = 43966
$ws.Range("B80").Formula = "=B79+C80"
$ws.Range("C80").Value = 9
$ws.Range("D80").Value = 1
$ws.Range("E80").Value = 7
$ws.Range("F80").Value = 5
$ws.Range("G80").Value = 31
$ws.Range("H80").Formula = "=G80+E80"
$ws.Range("I80").Formula = "=I79+J80"
$ws.Range("J80").Formula = "=K80+L80"
$ws.Range("K80").NumberFormat = "General"
$ws.Range("K80").Value = 0
$ws.Range("K80").NumberFormat = "@"
$ws.Range("L80").NumberFormat = "General"
$ws.Range("L80").Value = 0
$ws.Range("L80").NumberFormat = "@"

# Row 81 (43967)
$ws.Range("A81").Value = 43967
$ws.Range("B81").Formula = "=B80+C81"
$ws.Range("C81").Value = 3
$ws.Range("D81").Value = 0
$ws.Range("E81").Value = 7
$ws.Range("F81").Value = 5
$ws.Range("G81").Value = 31
$ws.Range("H81").Formula = "=G81+E81"
$ws.Range("I81").Formula = "=I80+J81"
$ws.Range("J81").Formula = "=K81+L81"
$ws.Range("K81").NumberFormat = "General"
$ws.Range("K81").Value = 0
$ws.Range("K81").NumberFormat = "@"
$ws.Range("L81").NumberFormat = "General"
$ws.Range("L81").Value = 0
$ws.Range("L81").NumberFormat = "@"

# Row 82 (43968)
$ws.Range("A82").Value = 43968
$ws.Range("B82").Formula = "=B81+C82"
$ws.Range("C82").Value = 0
$ws.Range("D82").Value = 0
$ws.Range("E82").Value = 7
$ws.Range("F82").Value = 5
$ws.Range("G82").Value = 31
$ws.Range("H82").Formula = "=G82+E82"
$ws.Range("I82").Formula = "=I81+J82"
$ws.Range("J82").Formula = "=K82+L82"
$ws.Range("K82").NumberFormat = "General"
$ws.Range("K82").Value = 0
$ws.Range("K82").NumberFormat = "@"
$ws.Range("L82").NumberFormat = "General"
$ws.Range("L82").Value = 0
$ws.Range("L82").NumberFormat = "@"

# Row 83 (43969)
$ws.Range("A83").Value = 43969
$ws.Range("B83").Formula = "=B82+C83"
$ws.Range("C83").Value = 1
$ws.Range("D83").Value = 1
$ws.Range("E83").Value = 6
$ws.Range("F83").Value = 4
$ws.Range("G83").Value = 30
$ws.Range("H83").Formula = "=G83+E83"
$ws.Range("I83").Formula = "=I82+J83"
$ws.Range("J83").Formula = "=K83+L83"
$ws.Range("K83").NumberFormat = "General"
$ws.Range("K83").Value = 1
$ws.Range("K83").NumberFormat = "@"
$ws.Range("L83").NumberFormat = "General"
$ws.Range("L83").Value = 1
$ws.Range("L83").NumberFormat = "@"

# Row 84 (43970)
$ws.Range("A84").Value = 43970
$ws.Range("B84").Formula = "=B83+C84"
$ws.Range("C84").Value = 1
$ws.Range("D84").Value = 1
$ws.Range("E84").Value = 6
$ws.Range("F84").Value = 4
$ws.Range("G84").Value = 29
$ws.Range("H84").Formula = "=G84+E84"
$ws.Range("I84").Formula = "=I83+J84"
$ws.Range("J84").Formula = "=K84+L84"
$ws.Range("K84").NumberFormat = "General"
$ws.Range("K84").Value = 0
$ws.Range("K84").NumberFormat = "@"
$ws.Range("L84").NumberFormat = "General"
$ws.Range("L84").Value = 0
$ws.Range("L84").NumberFormat = "@"

# Row 85 (43971)
$ws.Range("A85").Value = 43971
$ws.Range("B85").Formula = "=B84+C85"
$ws.Range("C85").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 5
$ws.Range("F85").Value = 4
$ws.Range("G85").Value = 28
$ws.Range("H85").Formula = "=G85+E85"
$ws.Range("I85").Formula = "=I84+J85"
$ws.Range("J85").Formula = "=K85+L85"
$ws.Range("K85").NumberFormat = "General"
$ws.Range("K85").Value = 0
$ws.Range("K85").NumberFormat = "@"
$ws.Range("L85").NumberFormat = "General"
$ws.Range("L85").Value = 0
$ws.Range("L85").NumberFormat = "@"

# Row 86 (43972)
$ws.Range("A86").Value = 43972
$ws.Range("B86").Formula = "=B85+C86"
$ws.Range("C86").Value = 1
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 5
$ws.Range("F86").Value = 4
$ws.Range("G86").Value = 28
$ws.Range("H86").Formula = "=G86+E86"
$ws.Range("I86").Formula = "=I85+J86"
$ws.Range("J86").Formula = "=K86+L86"
$ws.Range("K86").NumberFormat = "General"
$ws.Range("K86").Value = 0
$ws.Range("K86").NumberFormat = "@"
$ws.Range("L86").NumberFormat = "General"
$ws.Range("L86").Value = 0
$ws.Range("L86").NumberFormat = "@"

# Row 87 (43973)
$ws.Range("A87").Value = 43973
$ws.Range("B87").Formula = "=B86+C87"
$ws.Range("C87").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 4
$ws.Range("F87").Value = 3
$ws.Range("G87").Value = 28
$ws.Range("H87").Formula = "=G87+E87"
$ws.Range("I87").Formula = "=I86+J87"
$ws.Range("J87").Formula = "=K87+L87"
$ws.Range("K87").NumberFormat = "General"
$ws.Range("K87").Value = 0
$ws.Range("K87").NumberFormat = "@"
$ws.Range("L87").NumberFormat = "General"
$ws.Range("L87").Value = 0
$ws.Range("L87").NumberFormat = "@"

# Row 88 (43974)
$ws.Range("A88").Value = 43974
$ws.Range("B88").Formula = "=B87+C88"
$ws.Range("C88").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 4
$ws.Range("F88").Value = 3
$ws.Range("G88").Value = 28
$ws.Range("H88").Formula = "=G88+E88"
$ws.Range("I88").Formula = "=I87+J88"
$ws.Range("J88").Formula = "=K88+L88"
$ws.Range("K88").NumberFormat = "General"
$ws.Range("K88").Value = 0
$ws.Range("K88").NumberFormat = "@"
$ws.Range("L88").NumberFormat = "General"
$ws.Range("L88").Value = 0
$ws.Range("L88").NumberFormat = "@"

# Row 89 (43975)
$ws.Range("A89").Value = 43975
$ws.Range("B89").Formula = "=B88+C89"
$ws.Range("C89").Value = 1
$ws.Range("D89").Value = 0
$ws.Range("E89").Value = 4
$ws.Range("F89").Value = 3
$ws.Range("G89").Value = 28
$ws.Range("H89").Formula = "=G89+E89"
$ws.Range("I89").Formula = "=I88+J89"
$ws.Range("J89").Formula = "=K89+L89"
$ws.Range("K89").NumberFormat = "General"
$ws.Range("K89").Value = 0
$ws.Range("K89").NumberFormat = "@"
$ws.Range("L89").NumberFormat = "General"
$ws.Range("L89").Value = 0
$ws.Range("L89").NumberFormat = "@"

# Row 90 (new placeholder / "today" row - B/C left blank, like old row 85)
$ws.Range("A90").Value = 43976
$ws.Range("D90").Value = 0
$ws.Range("E90").Value = 4
$ws.Range("F90").Value = 3
$ws.Range("G90").Value = 28
$ws.Range("H90").Formula = "=G90+E90"
$ws.Range("I90").Formula = "=I89+J90"
$ws.Range("J90").Formula = "=K90+L90"
$ws.Range("K90").NumberFormat = "General"
$ws.Range("K90").Value = 0
$ws.Range("K90").NumberFormat = "@"
$ws.Range("L90").NumberFormat = "General"
$ws.Range("L90").Value = 0
$ws.Range("L90").NumberFormat = "@"

# --- 5. Move the view/selection the way the author left it -----------------
$ws.Range("A76").Select()
